$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New formulas for rows 33-34 (columns E and F were previously empty) ---
$ws.Range("E33").NumberFormat = "0%"
$ws.Range("E33").Formula = '=AVERAGE(C33:$C$41)/$C$43'
$ws.Range("F33").NumberFormat = "0%"
$ws.Range("F33").Formula = '=C33/$C$43'

$ws.Range("E34").NumberFormat = "0%"
$ws.Range("E34").Formula = '=AVERAGE(C34:$C$41)/$C$43'
$ws.Range("F34").NumberFormat = "0%"
$ws.Range("F34").Formula = '=C34/$C$43'

# --- Fill in previously-empty F formulas for rows 35-41 ---
$ws.Range("F35").Formula = '=C35/$C$43'
$ws.Range("F36").Formula = '=C36/$C$43'
$ws.Range("F37").Formula = '=C37/$C$43'
$ws.Range("F38").Formula = '=C38/$C$43'
$ws.Range("F39").Formula = '=C39/$C$43'
$ws.Range("F40").Formula = '=C40/$C$43'
$ws.Range("F41").Formula = '=C41/$C$43'

# --- Update selection to F34:F35 (also clears the old scrolled topLeftCell) ---
$ws.Range("F34:F35").Select()
